$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1799
$ws.Range("I40").Value = 999
$ws.Range("J40").Value = 2332.3333
$ws.Range("K40").Value = 999
$ws.Range("L40").Value = 2332.3333
$ws.Range("M40").Value = -824
$ws.Range("N40").Value = -2682.3333
$ws.Range("H69").Value = 10935.8
$ws.Range("I69").Value = 9920
$ws.Range("K69").Value = 29760
$ws.Range("M69").Value = -28886
$ws.Range("H72").Value = 10935.8
$ws.Range("I72").Value = 9920
$ws.Range("K72").Value = 89280
$ws.Range("M72").Value = -84912
$ws.Range("H106").Value = 19999.75
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H132").Value = 1559.375
$ws.Range("I132").Value = 1596.7391
$ws.Range("K132").Value = 4790.2173
$ws.Range("M132").Value = -2260.2173
$ws.Range("H137").Value = 4189.457
$ws.Range("I137").Value = 1363.9354
$ws.Range("J137").Value = 26087.25
$ws.Range("K137").Value = 4091.8062
$ws.Range("L137").Value = 78261.75
$ws.Range("M137").Value = -1541.8062
$ws.Range("N137").Value = -83361.75
$ws.Range("H138").Value = 4341.2856
$ws.Range("I138").Value = 4194.76
$ws.Range("J138").Value = 4493.9165
$ws.Range("K138").Value = 12584.28
$ws.Range("L138").Value = 13481.7495
$ws.Range("M138").Value = -7444.280000000001
$ws.Range("N138").Value = -23761.7495
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 26999
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H32").Value = 6025.315
$ws.Range("J32").Value = 25519.8
$ws.Range("L32").Value = 25519.8
$ws.Range("N32").Value = -26093.8
$ws.Range("H45").Value = 4544
$ws.Range("J45").Value = 5266.5
$ws.Range("L45").Value = 5266.5
$ws.Range("N45").Value = -6020.5
$ws.Range("H116").Value = 26999
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H132").Value = 4647.1816
$ws.Range("I132").Value = 3645.25
$ws.Range("K132").Value = 10935.75
$ws.Range("M132").Value = -8405.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 26999
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H35").Value = 37747
$ws.Range("J35").Value = 39999
$ws.Range("L35").Value = 39999
$ws.Range("N35").Value = -40619
$ws.Range("H80").Value = 851.8570999999999
$ws.Range("J80").Value = 963.1
$ws.Range("L80").Value = 963.1
$ws.Range("N80").Value = -2959.1
$ws.Range("H83").Value = 851.8570999999999
$ws.Range("J83").Value = 963.1
$ws.Range("L83").Value = 4815.5
$ws.Range("N83").Value = -14799.5
$ws.Range("H99").Value = 6543.45
$ws.Range("I99").Value = 6730
$ws.Range("K99").Value = 6730
$ws.Range("M99").Value = -5232
$ws.Range("H132").Value = 76708.664
$ws.Range("J132").Value = 76708.5
$ws.Range("L132").Value = 76708.5
$ws.Range("N132").Value = -86828.5
$ws.Range("H140").Value = 69999.75
$ws.Range("J140").Value = 68333
$ws.Range("L140").Value = 68333
$ws.Range("N140").Value = -78693
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2661.7144
$ws.Range("I31").Value = 2105.3333
$ws.Range("K31").Value = 2105.3333
$ws.Range("M31").Value = -1810.3333
$ws.Range("H34").Value = 2661.7144
$ws.Range("I34").Value = 2105.3333
$ws.Range("K34").Value = 2105.3333
$ws.Range("M34").Value = -1903.3333
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H86").Value = 13908.852
$ws.Range("I86").Value = 14260.375
$ws.Range("K86").Value = 14260.375
$ws.Range("M86").Value = -13137.375
$ws.Range("H89").Value = 13908.852
$ws.Range("I89").Value = 14260.375
$ws.Range("K89").Value = 71301.875
$ws.Range("M89").Value = -65685.875
$ws.Range("H105").Value = 33845
$ws.Range("I105").Value = 100000
$ws.Range("K105").Value = 100000
$ws.Range("M105").Value = -98253
$ws.Range("H107").Value = 742.61536
$ws.Range("J107").Value = 960.2222
$ws.Range("L107").Value = 960.2222
$ws.Range("N107").Value = -4800.2222
$ws.Range("H122").Value = 27917.55
$ws.Range("I122").Value = 3314.8235
$ws.Range("J122").Value = 167333
$ws.Range("K122").Value = 9944.470499999999
$ws.Range("L122").Value = 501999
$ws.Range("M122").Value = -7494.470499999999
$ws.Range("N122").Value = -506899
$ws.Range("H132").Value = 2897.0657
$ws.Range("I132").Value = 2705.9272
$ws.Range("J132").Value = 4649.1665
$ws.Range("K132").Value = 8117.7816
$ws.Range("L132").Value = 13947.4995
$ws.Range("M132").Value = -5587.7816
$ws.Range("N132").Value = -19007.4995
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 68234.12
$ws.Range("J37").Value = 68234.12
$ws.Range("L37").Value = 204702.36
$ws.Range("N37").Value = -204926.36
$ws.Range("H86").Value = 1299.6666
$ws.Range("I86").Value = 1449.5
$ws.Range("K86").Value = 4348.5
$ws.Range("M86").Value = -3162.5
$ws.Range("H89").Value = 1299.6666
$ws.Range("I89").Value = 1449.5
$ws.Range("K89").Value = 13045.5
$ws.Range("M89").Value = -7117.5
$ws.Range("H107").Value = 468
$ws.Range("I107").Value = 199.625
$ws.Range("J107").Value = 1004.75
$ws.Range("K107").Value = 598.875
$ws.Range("L107").Value = 3014.25
$ws.Range("M107").Value = 1321.125
$ws.Range("N107").Value = -6854.25
$ws.Range("H113").Value = 1620.4839
$ws.Range("I113").Value = 507.5
$ws.Range("K113").Value = 1522.5
$ws.Range("M113").Value = 647.5
$ws.Range("H132").Value = 3374.1667
$ws.Range("J132").Value = 3062.5
$ws.Range("L132").Value = 27562.5
$ws.Range("N132").Value = -32622.5
$ws.Range("H137").Value = 5219.778
$ws.Range("J137").Value = 7797.8
$ws.Range("L137").Value = 23393.4
$ws.Range("N137").Value = -33593.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1499.8695
$ws.Range("J113").Value = 1364.0714
$ws.Range("L113").Value = 1364.0714
$ws.Range("N113").Value = -5704.0714
$ws.Range("H132").Value = 15941
$ws.Range("I132").Value = 15941
$ws.Range("K132").Value = 47823
$ws.Range("M132").Value = -45293
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3089.2
$ws.Range("J22").Value = 3428.15
$ws.Range("L22").Value = 3428.15
$ws.Range("N22").Value = -4018.15
$ws.Range("H27").Value = 3089.2
$ws.Range("J27").Value = 3428.15
$ws.Range("L27").Value = 3428.15
$ws.Range("N27").Value = -3642.15
$ws.Range("H46").Value = 4859.6
$ws.Range("I46").Value = 1050.5
$ws.Range("J46").Value = 5445.615
$ws.Range("K46").Value = 1050.5
$ws.Range("L46").Value = 5445.615
$ws.Range("M46").Value = -862.5
$ws.Range("N46").Value = -5821.615
$ws.Range("H122").Value = 5566.3184
$ws.Range("J122").Value = 6482.5
$ws.Range("L122").Value = 19447.5
$ws.Range("N122").Value = -24347.5
$ws.Range("H133").Value = 88990.2
$ws.Range("J133").Value = 88990.2
$ws.Range("L133").Value = 88990.2
$ws.Range("N133").Value = -94050.2
$ws.Range("H136").Value = 3170
$ws.Range("J136").Value = 3219
$ws.Range("L136").Value = 9657
$ws.Range("N136").Value = -14757
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 25080.9
$ws.Range("I122").Value = 2359.6155
$ws.Range("K122").Value = 7078.8465
$ws.Range("M122").Value = -4628.8465
$ws.Range("H125").Value = 59374.125
$ws.Range("J125").Value = 59374.125
$ws.Range("L125").Value = 59374.125
$ws.Range("N125").Value = -69214.125
